$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume cells stay plain text (matches source data which is
# stored as inline strings, e.g. "6.60" must not become the number 6.6).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '64.969.72'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '3.173.27'
$ws.Range("E3").Value = '  +3.64%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '577.02'
$ws.Range("E5").Value = '  +2.95%  '
$ws.Range("D6").Value = '150.67'
$ws.Range("E6").Value = '  +4.96%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.170.68'
$ws.Range("E8").Value = '  +3.57%  '
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  +3.34%  '
$ws.Range("E10").Value = '  +5.11%  '
$ws.Range("D11").Value = '6.19'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").Value = '0.504'
$ws.Range("E12").Value = '  +1.94%  '
$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").Value = '  +16.93%  '
$ws.Range("D14").Value = '37.92'
$ws.Range("E14").Value = '  +6.10%  '
$ws.Range("D15").Value = '3.698.93'
$ws.Range("E15").Value = '  +3.78%  '
$ws.Range("D16").Value = '65.068.25'
$ws.Range("E16").Value = '  +1.58%  '
$ws.Range("D17").Value = '3.177.16'
$ws.Range("E17").Value = '  +3.73%  '
$ws.Range("D18").Value = '7.16'
$ws.Range("E18").Value = '  +5.29%  '
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").Value = '511.28'
$ws.Range("E20").Value = '  +6.93%  '
$ws.Range("D21").Value = '14.84'
$ws.Range("E21").Value = '  +6.09%  '
$ws.Range("D22").Value = '0.732'
$ws.Range("E22").Value = '  +6.80%  '
$ws.Range("D23").Value = '15.27'
$ws.Range("E23").Value = '  +6.59%  '
$ws.Range("D24").Value = '7.81'
$ws.Range("E24").Value = '  +3.09%  '
$ws.Range("D25").Value = '85.27'
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").Value = '9.02'
$ws.Range("E27").Value = '  +11.37%  '
$ws.Range("D28").Value = '2.93'
$ws.Range("E28").Value = '  +4.32%  '
$ws.Range("D29").Value = '2.18'
$ws.Range("E29").Value = '  +6.71%  '
$ws.Range("D30").Value = '27.98'
$ws.Range("E30").Value = '  +6.19%  '
$ws.Range("D31").Value = '2.79'
$ws.Range("E31").Value = '  +13.43%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("E33").Value = '  +5.48%  '
$ws.Range("D34").Value = '6.28'
$ws.Range("E34").Value = '  +8.52%  '
$ws.Range("D35").Value = '6.60'
$ws.Range("E35").Value = '  +5.83%  '
$ws.Range("D36").Value = '55.60'
$ws.Range("E36").Value = '  +1.66%  '
$ws.Range("D37").Value = '0.0897'
$ws.Range("E37").Value = '  +10.01%  '
$ws.Range("D38").Value = '475.48'
$ws.Range("E38").Value = '  +5.21%  '
$ws.Range("D39").Value = '3.13'
$ws.Range("E39").Value = '  +10.38%  '
$ws.Range("D40").Value = '0.0421'
$ws.Range("E40").Value = '  +2.20%  '
$ws.Range("D41").Value = '8.64'
$ws.Range("E41").Value = '  +4.20%  '
$ws.Range("D42").Value = '3.057.93'
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("E43").Value = '  +2.36%  '
$ws.Range("D44").Value = '0.287'
$ws.Range("E44").Value = '  +7.46%  '
$ws.Range("D45").Value = '2.39'
$ws.Range("E45").Value = '  +7.34%  '
$ws.Range("D46").Value = '29.07'
$ws.Range("E46").Value = '  +4.38%  '
$ws.Range("D47").Value = '0.0₃0611'
$ws.Range("E47").Value = '  +18.12%  '
$ws.Range("D49").Value = '0.115'
$ws.Range("E49").Value = '  +1.61%  '
$ws.Range("D50").Value = '2.26'
$ws.Range("E50").Value = '  +7.61%  '
$ws.Range("D51").Value = '120.23'
$ws.Range("E51").Value = '  +0.72%  '
